$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 56, shifting existing rows 56-92 down to 57-93.
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new price observation.
$ws.Range("A56").Value = 11
$ws.Range("B56").Value = "Vega Monumental Concepción"
$ws.Range("C56").Value = "Bíobío"
$ws.Range("D56").Value = "2022-09-09"
$ws.Range("E56").Value = 8
$ws.Range("F56").Value = 100112012
$ws.Range("G56").Value = "Espinaca"
$ws.Range("H56").Value = "Sin especificar"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 50
$ws.Range("K56").Value = 6500
$ws.Range("L56").Value = 7000
$ws.Range("M56").Value = 6700
$ws.Range("N56").Value = "$/cuna 10 kilos"
$ws.Range("O56").Value = "Región Metropolitana"
$ws.Range("P56").Value = 670
$ws.Range("Q56").Value = 10
$ws.Range("R56").Value = "Hortaliza"
